$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 6
$ws.Range("H6").Value = 1200
$ws.Range("I6").Value = 400
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 1200
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = -1088
$ws.Range("N6").Value = -6224
# Row 48
$ws.Range("H48").Value = 9186.875
$ws.Range("I48").Value = 7749.5
$ws.Range("J48").Value = 10624.25
$ws.Range("K48").Value = 23248.5
$ws.Range("L48").Value = 31872.75
$ws.Range("M48").Value = -22956.5
$ws.Range("N48").Value = -32456.75
# Row 49
$ws.Range("H49").Value = 4173.8
$ws.Range("I49").Value = 425
$ws.Range("K49").Value = 1275
$ws.Range("M49").Value = -1139
# Row 56
$ws.Range("H56").Value = 9186.875
$ws.Range("I56").Value = 7749.5
$ws.Range("J56").Value = 10624.25
$ws.Range("K56").Value = 23248.5
$ws.Range("L56").Value = 31872.75
$ws.Range("M56").Value = -22714.5
$ws.Range("N56").Value = -32940.75
# Row 88
$ws.Range("H88").Value = 13414.167
$ws.Range("I88").Value = 3247
$ws.Range("J88").Value = 18497.75
$ws.Range("K88").Value = 3247
$ws.Range("L88").Value = 18497.75
$ws.Range("M88").Value = -2841
$ws.Range("N88").Value = -19309.75
# Row 91
$ws.Range("H91").Value = 13414.167
$ws.Range("I91").Value = 3247
$ws.Range("J91").Value = 18497.75
$ws.Range("K91").Value = 3247
$ws.Range("L91").Value = 18497.75
$ws.Range("M91").Value = -1843
$ws.Range("N91").Value = -21305.75
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# Row 141
$ws.Range("H141").Value = 5159.7334
$ws.Range("I141").Value = 732.6667
$ws.Range("K141").Value = 2198.0001
$ws.Range("M141").Value = 2981.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4503.4316
$ws.Range("I32").Value = 3933.4187
$ws.Range("K32").Value = 3933.4187
$ws.Range("M32").Value = -3646.4187
# Row 132
$ws.Range("H132").Value = 2884.7778
$ws.Range("I132").Value = 2311.8635
$ws.Range("J132").Value = 5405.6
$ws.Range("K132").Value = 6935.5905
$ws.Range("L132").Value = 16216.8
$ws.Range("M132").Value = -4405.5905
$ws.Range("N132").Value = -21276.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1956.5264
$ws.Range("I134").Value = 1199.4375
$ws.Range("J134").Value = 5994.3335
$ws.Range("K134").Value = 3598.3125
$ws.Range("L134").Value = 17983.0005
$ws.Range("M134").Value = -1063.3125
$ws.Range("N134").Value = -23053.0005
# Row 140
$ws.Range("H140").Value = 68147.336
$ws.Range("J140").Value = 68147.336
$ws.Range("L140").Value = 68147.336
$ws.Range("N140").Value = -78507.336

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 30891.486
$ws.Range("I31").Value = 3487.8
$ws.Range("J31").Value = 59737.473
$ws.Range("K31").Value = 3487.8
$ws.Range("L31").Value = 59737.473
$ws.Range("M31").Value = -3192.8
$ws.Range("N31").Value = -60327.473
# Row 34
$ws.Range("H34").Value = 30891.486
$ws.Range("I34").Value = 3487.8
$ws.Range("J34").Value = 59737.473
$ws.Range("K34").Value = 3487.8
$ws.Range("L34").Value = 59737.473
$ws.Range("M34").Value = -3285.8
$ws.Range("N34").Value = -60141.473
# Row 50
$ws.Range("H50").Value = 58950
$ws.Range("J50").Value = 58950
$ws.Range("L50").Value = 58950
$ws.Range("N50").Value = -60200
# Row 58
$ws.Range("H58").Value = 4973.778
$ws.Range("I58").Value = 2493
$ws.Range("J58").Value = 8074.75
$ws.Range("K58").Value = 2493
$ws.Range("L58").Value = 8074.75
$ws.Range("M58").Value = -2290
$ws.Range("N58").Value = -8480.75
# Row 94
$ws.Range("H94").Value = 2227.9092
$ws.Range("J94").Value = 2526.2856
$ws.Range("L94").Value = 2526.2856
$ws.Range("N94").Value = -3428.2856
# Row 134
$ws.Range("H134").Value = 5158.1113
$ws.Range("I134").Value = 3499.8333
$ws.Range("K134").Value = 10499.4999
$ws.Range("M134").Value = -7964.499899999999
# Row 136
$ws.Range("H136").Value = 4973.778
$ws.Range("I136").Value = 2493
$ws.Range("J136").Value = 8074.75
$ws.Range("K136").Value = 7479
$ws.Range("L136").Value = 24224.25
$ws.Range("M136").Value = -4929
$ws.Range("N136").Value = -29324.25
# Row 141
$ws.Range("H141").Value = 140325.64
$ws.Range("J141").Value = 140325.64
$ws.Range("L141").Value = 140325.64
$ws.Range("N141").Value = -150685.64

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 138
$ws.Range("H138").Value = 29717.8
$ws.Range("I138").Value = 38696.332
$ws.Range("J138").Value = 16250
$ws.Range("K138").Value = 116088.996
$ws.Range("L138").Value = 48750
$ws.Range("M138").Value = -110948.996
$ws.Range("N138").Value = -59030

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
# Row 70
$ws.Range("H70").Value = 20247.066
$ws.Range("I70").Value = 6079
$ws.Range("K70").Value = 6079
$ws.Range("M70").Value = -5809
# Row 73
$ws.Range("H73").Value = 20247.066
$ws.Range("I73").Value = 6079
$ws.Range("K73").Value = 6079
$ws.Range("M73").Value = -5143
# Row 126
$ws.Range("H126").Value = 3978.45
$ws.Range("I126").Value = 3185.5881
$ws.Range("K126").Value = 9556.764299999999
$ws.Range("M126").Value = -7086.764299999999
# Row 132
$ws.Range("H132").Value = 50608.184
$ws.Range("I132").Value = 79983.766
$ws.Range("J132").Value = 8176.778
$ws.Range("K132").Value = 239951.298
$ws.Range("L132").Value = 24530.334
$ws.Range("M132").Value = -237421.298
$ws.Range("N132").Value = -29590.334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6969.8
$ws.Range("I7").Value = 2853.111
$ws.Range("J7").Value = 10338
$ws.Range("K7").Value = 2853.111
$ws.Range("L7").Value = 10338
$ws.Range("M7").Value = -2741.111
$ws.Range("N7").Value = -10562
# Row 126
$ws.Range("H126").Value = 6969.8
$ws.Range("I126").Value = 2853.111
$ws.Range("J126").Value = 10338
$ws.Range("K126").Value = 8559.332999999999
$ws.Range("L126").Value = 31014
$ws.Range("M126").Value = -6089.332999999999
$ws.Range("N126").Value = -35954

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 2074.818
$ws.Range("I136").Value = 1552.1333
$ws.Range("K136").Value = 4656.3999
$ws.Range("M136").Value = -2106.3999
